# ExcelManager: fixed addProductToInventory and reduceProductQuantity
# methods are reflected here by rebuilding the inventory table -
# the obsolete "S/No." column is dropped, the existing rows are
# refreshed with correct product data, and a new row is appended
# for a product (Beer) added via addProductToInventory.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "S/No." column entirely - Name/Price/Quantity/Category
# shift left from B:E into A:D.
$ws.Columns.Item(1).Delete()

# Row 1 (header) is unchanged content-wise: Name, Price, Quantity, Category

# Row 2: Rice - price & quantity corrected, category re-affirmed.
$ws.Range("A2").Value = "Rice"
$ws.Range("B2").Value = 50000
$ws.Range("C2").Value = 250
$ws.Range("D2").Value = "Groceries"

# Row 3: was Beer/Drinks - replaced with Jeans/Fashion.
$ws.Range("A3").Value = "Jeans"
$ws.Range("B3").Value = 12000
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "Fashion"

# Row 4: was Soap/Toiletries - replaced with LGTV/Electronics.
$ws.Range("A4").Value = "LGTV"
$ws.Range("B4").Value = 270000
$ws.Range("C4").Value = 80
$ws.Range("D4").Value = "Electronics"

# Row 5: newly added product via addProductToInventory - Beer.
$ws.Range("A5").Value = "Beer"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "Groceries"

# Leave the selection where the editor's cursor landed after the edit.
[void]$ws.Range("D7").Select()
